$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row text was relabelled ("素"/"香" placeholders -> real column headers)
$ws.Range("B1").Value = "點餐者"
$ws.Range("C1").Value = "時間"

# The old 4th column ("大雞排") is dropped entirely from the sheet
$ws.Range("D1").ClearContents()

# New explicit column widths for the two renamed columns
$ws.Columns("B").ColumnWidth = 13.95
$ws.Columns("C").ColumnWidth = 20.8

# Selection moved to F3
$ws.Range("F3").Select()
